# Update view-count column F for several rows across sheets
# (gh-pages data refresh: 'Update gh-pages to output generated at 456a3b4')
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1861
$ws.Range("F5").Value = 3226
$ws.Range("F6").Value = 200
$ws.Range("F7").Value = 4714
$ws.Range("F8").Value = 446
$ws.Range("F9").Value = 265
$ws.Range("F10").Value = 157
$ws.Range("F13").Value = 13
$ws.Range("F14").Value = 7
$ws.Range("F15").Value = 641
$ws.Range("F16").Value = 279
$ws.Range("F18").Value = 85
$ws.Range("F19").Value = 137
$ws.Range("F21").Value = 4675
$ws.Range("F25").Value = 5823
$ws.Range("F27").Value = 1176
$ws.Range("F29").Value = 650
$ws.Range("F30").Value = 4406
$ws.Range("F32").Value = 71
$ws.Range("F33").Value = 114
$ws.Range("F34").Value = 795
$ws.Range("F35").Value = 56
$ws.Range("F36").Value = 725
$ws.Range("F37").Value = 733

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 32

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 1861
$ws.Range("F9").Value = 3226
$ws.Range("F10").Value = 200
$ws.Range("F11").Value = 4714
$ws.Range("F12").Value = 446
$ws.Range("F13").Value = 265
$ws.Range("F14").Value = 157
$ws.Range("F17").Value = 13
$ws.Range("F18").Value = 7
$ws.Range("F19").Value = 641
$ws.Range("F20").Value = 279
$ws.Range("F22").Value = 32
$ws.Range("F23").Value = 85
$ws.Range("F24").Value = 137
$ws.Range("F26").Value = 4675
$ws.Range("F30").Value = 5823
$ws.Range("F32").Value = 1176
$ws.Range("F34").Value = 650
$ws.Range("F35").Value = 4406
$ws.Range("F38").Value = 71
$ws.Range("F39").Value = 114
$ws.Range("F40").Value = 795
$ws.Range("F41").Value = 56
$ws.Range("F42").Value = 725
$ws.Range("F43").Value = 733
